$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.362.73'
$ws.Range('E2').Value = '  -0.49%  '
$ws.Range('D3').Value = '1.846.48'
$ws.Range('E3').Value = '  -0.24%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9989'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '240.28'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.67%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.6304'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.02%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.07542'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.2953'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.89%  '
$ws.Range('E10').Value = '  -0.07%  '
$ws.Range('E11').Value = '  -0.26%  '
$ws.Range('D12').Value = '1.871.37'
$ws.Range('E12').Value = '  -0.11%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.6826'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.38%  '
$ws.Range('E15').Value = '  +2.69%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '82.86'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.01%  '
$ws.Range('B17').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C17').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D17').Value = '2.118.27'
$ws.Range('E17').Value = '  -1.37%  '
$ws.Range('B18').Value = 'Uniswap'
$ws.Range('C18').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '6.127'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.91%  '
$ws.Range('B19').Value = 'WrappedBTC'
$ws.Range('C19').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D19').Value = '29.393.41'
$ws.Range('E19').Value = '  -0.56%  '
$ws.Range('B20').Value = 'BitcoinCash'
$ws.Range('C20').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '227.22'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -3.06%  '
$ws.Range('B21').Value = 'Avalanche'
$ws.Range('C21').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.44'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.55%  '
$ws.Range('B22').Value = 'Dai'
$ws.Range('C22').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.9998'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.06%  '
$ws.Range('B23').Value = 'Chainlink'
$ws.Range('C23').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.538'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.43%  '
$ws.Range('E24').Value = '  +0.11%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '157.09'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.51%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.1393'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.05%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.352'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.32%  '
$ws.Range('E28').Value = '  -0.54%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.464'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.01%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.05677'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -3.39%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.252'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.02%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.119'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.26%  '
$ws.Range('E33').Value = '  -0.80%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.845'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.85%  '
$ws.Range('E35').Value = '  -1.43%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7121'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.24%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.594'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.23%  '
$ws.Range('E38').Value = '  +1.34%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01814'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.63%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.780'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.43%  '
$ws.Range('B41').Value = 'FraxShare'
$ws.Range('C41').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.211'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.73%  '
$ws.Range('B42').Value = 'TrustWalletToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.9095'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.06%  '
$ws.Range('E43').Value = '  +0.08%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '101.21'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.71%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '66.15'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.75%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '7.066'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -4.41%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.4037'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.25%  '
$ws.Range('B48').Value = 'BabyDogeCoin'
$ws.Range('C48').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.00000000117'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.60%  '
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '9.059'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.24%  '
$ws.Range('B50').Value = 'RenderToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.679'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.39%  '
$ws.Range('B51').Value = 'Algorand'
$ws.Range('C51').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.1123'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.40%  '
